$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 54.86839566666666
$ws.Range("H2").Value = 164.605187
$ws.Range("I2").Value = 0.6170939026906647
$ws.Range("J2").Value = 0.6170939026906647
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 7994.366107298641
$ws.Range("R2").Value = 71949.29496568778
$ws.Range("S2").Value = 0.1768552711081317
$ws.Range("T2").Value = 0.1768552711081317

$ws.Range("G3").Value = 54.86839566666666
$ws.Range("H3").Value = 164.605187
$ws.Range("I3").Value = 0.6170939026906647
$ws.Range("J3").Value = 0.6170939026906647
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 9261.768874330355
$ws.Range("R3").Value = 83355.91986897321
$ws.Range("S3").Value = 0.2048933740619044
$ws.Range("T3").Value = 0.2048933740619044

$ws.Range("G4").Value = 54.86839566666666
$ws.Range("H4").Value = 164.605187
$ws.Range("I4").Value = 0.6170939026906647
$ws.Range("J4").Value = 0.6170939026906647
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 7030.076549499197
$ws.Range("R4").Value = 63270.68894549277
$ws.Range("S4").Value = 0.155522786595612
$ws.Range("T4").Value = 0.155522786595612

$ws.Range("G5").Value = 54.86839566666666
$ws.Range("H5").Value = 164.605187
$ws.Range("I5").Value = 0.6170939026906647
$ws.Range("J5").Value = 0.6170939026906647
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 3608.204902038924
$ws.Range("R5").Value = 32473.84411835032
$ws.Range("S5").Value = 0.07982247092501658
$ws.Range("T5").Value = 0.07982247092501658

$ws.Range("I6").Value = 0.06720170646055251
$ws.Range("J6").Value = 0.0672017064605525
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 870.5888068872671
$ws.Range("R6").Value = 7835.299261985404
$ws.Range("S6").Value = 0.01925959074168291
$ws.Range("T6").Value = 0.01925959074168291

$ws.Range("I7").Value = 0.06720170646055251
$ws.Range("J7").Value = 0.0672017064605525
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("S7").Value = 0.02231294835256616
$ws.Range("T7").Value = 0.02231294835256615

$ws.Range("I8").Value = 0.06720170646055251
$ws.Range("J8").Value = 0.0672017064605525
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 765.5773920545104
$ws.Range("R8").Value = 6890.196528490593
$ws.Range("S8").Value = 0.01693647693998448
$ws.Range("T8").Value = 0.01693647693998447

$ws.Range("I9").Value = 0.06720170646055251
$ws.Range("J9").Value = 0.0672017064605525
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 392.9345689838957
$ws.Range("R9").Value = 3536.411120855061
$ws.Range("S9").Value = 0.008692690426318973
$ws.Range("T9").Value = 0.008692690426318971

$ws.Range("G10").Value = 4.832157666666667
$ws.Range("H10").Value = 14.496473
$ws.Range("I10").Value = 0.05434631351453007
$ws.Range("J10").Value = 0.05434631351453007
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 704.048970380076
$ws.Range("R10").Value = 6336.440733420684
$ws.Range("S10").Value = 0.01557531514803791
$ws.Range("T10").Value = 0.01557531514803791

$ws.Range("G11").Value = 4.832157666666667
$ws.Range("H11").Value = 14.496473
$ws.Range("I11").Value = 0.05434631351453007
$ws.Range("J11").Value = 0.05434631351453007
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 815.6667773717874
$ws.Range("R11").Value = 7341.000996346086
$ws.Range("S11").Value = 0.01804457878333627
$ws.Range("T11").Value = 0.01804457878333627

$ws.Range("G12").Value = 4.832157666666667
$ws.Range("H12").Value = 14.496473
$ws.Range("I12").Value = 0.05434631351453007
$ws.Range("J12").Value = 0.05434631351453007
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 619.1257805730526
$ws.Range("R12").Value = 5572.132025157473
$ws.Range("S12").Value = 0.01369660287052832
$ws.Range("T12").Value = 0.01369660287052832

$ws.Range("G13").Value = 4.832157666666667
$ws.Range("H13").Value = 14.496473
$ws.Range("I13").Value = 0.05434631351453007
$ws.Range("J13").Value = 0.05434631351453007
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 317.7679020581223
$ws.Range("R13").Value = 2859.911118523101
$ws.Range("S13").Value = 0.007029816712627581
$ws.Range("T13").Value = 0.007029816712627581

$ws.Range("G14").Value = 23.238438
$ws.Range("H14").Value = 69.71531400000001
$ws.Range("I14").Value = 0.2613580773342528
$ws.Range("J14").Value = 0.2613580773342528
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 3385.857721490165
$ws.Range("R14").Value = 30472.71949341149
$ws.Range("S14").Value = 0.07490359801273172
$ws.Range("T14").Value = 0.07490359801273172

$ws.Range("G15").Value = 23.238438
$ws.Range("H15").Value = 69.71531400000001
$ws.Range("I15").Value = 0.2613580773342528
$ws.Range("J15").Value = 0.2613580773342528
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 3922.641424837769
$ws.Range("R15").Value = 35303.77282353992
$ws.Range("S15").Value = 0.08677858923877729
$ws.Range("T15").Value = 0.08677858923877729

$ws.Range("G16").Value = 23.238438
$ws.Range("H16").Value = 69.71531400000001
$ws.Range("I16").Value = 0.2613580773342528
$ws.Range("J16").Value = 0.2613580773342528
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 2977.451701399745
$ws.Range("R16").Value = 26797.0653125977
$ws.Range("S16").Value = 0.06586864058948569
$ws.Range("T16").Value = 0.06586864058948569

$ws.Range("G17").Value = 23.238438
$ws.Range("H17").Value = 69.71531400000001
$ws.Range("I17").Value = 0.2613580773342528
$ws.Range("J17").Value = 0.2613580773342528
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 1528.184757154602
$ws.Range("R17").Value = 13753.66281439142
$ws.Range("S17").Value = 0.03380724949325808
$ws.Range("T17").Value = 0.03380724949325808

